# Update cryptocurrency price/volume data (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.018.92'
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("D3").Value = '1.555.42'
$ws.Range("E3").Value = '  -0.80%  '
$c = $ws.Range("D4")
$c.Value = "'0.9996"
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.14%  '
$c = $ws.Range("D5")
$c.Value = "'0.9998"
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.13%  '
$c = $ws.Range("D6")
$c.Value = "'286.97"
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.09%  '
$c = $ws.Range("D7")
$c.Value = "'0.3825"
$c.ClearFormats()
$ws.Range("E7").Value = '  +3.60%  '
$c = $ws.Range("D8")
$c.Value = "'0.3239"
$c.ClearFormats()
$ws.Range("E8").Value = '  -1.71%  '
$c = $ws.Range("D9")
$c.Value = "'41.37"
$c.ClearFormats()
$ws.Range("E9").Value = '  -12.31%  '
$c = $ws.Range("D10")
$c.Value = "'1.122"
$c.ClearFormats()
$ws.Range("E10").Value = '  -2.61%  '
$c = $ws.Range("D11")
$c.Value = "'0.07316"
$c.ClearFormats()
$ws.Range("E11").Value = '  -1.66%  '
$c = $ws.Range("D12")
$c.Value = "'0.9997"
$c.ClearFormats()
$ws.Range("E12").Value = '  -0.15%  '
$c = $ws.Range("D13")
$c.Value = "'19.38"
$c.ClearFormats()
$ws.Range("E13").Value = '  -6.18%  '
$c = $ws.Range("D14")
$c.Value = "'5.722"
$c.ClearFormats()
$ws.Range("E14").Value = '  -2.98%  '
$c = $ws.Range("D15")
$c.Value = "'6.816"
$c.ClearFormats()
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '1.555.09'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("E17").Value = '  -1.39%  '
$c = $ws.Range("D18")
$c.Value = "'0.06623"
$c.ClearFormats()
$ws.Range("E18").Value = '  -1.23%  '
$c = $ws.Range("D19")
$c.Value = "'85.23"
$c.ClearFormats()
$ws.Range("E19").Value = '  -1.63%  '
$c = $ws.Range("D20")
$c.Value = "'6.407"
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.79%  '
$c = $ws.Range("D21")
$c.Value = "'0.9997"
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.13%  '
$c = $ws.Range("D22")
$c.Value = "'15.95"
$c.ClearFormats()
$ws.Range("E22").Value = '  -2.69%  '
$c = $ws.Range("D23")
$c.Value = "'11.45"
$c.ClearFormats()
$ws.Range("E23").Value = '  -3.71%  '
$ws.Range("D24").Value = '22.031.27'
$ws.Range("E24").Value = '  -1.65%  '
$c = $ws.Range("D25")
$c.Value = "'2.294"
$c.ClearFormats()
$ws.Range("E25").Value = '  -2.52%  '
$c = $ws.Range("D26")
$c.Value = "'2.527"
$c.ClearFormats()
$ws.Range("E26").Value = '  -2.69%  '
$c = $ws.Range("D27")
$c.Value = "'148.91"
$c.ClearFormats()
$ws.Range("E27").Value = '  -1.40%  '
$c = $ws.Range("D28")
$c.Value = "'18.82"
$c.ClearFormats()
$ws.Range("E28").Value = '  -3.43%  '
$c = $ws.Range("D29")
$c.Value = "'4.856"
$c.ClearFormats()
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("D30").Value = '1.729.01'
$ws.Range("E30").Value = '  -0.94%  '
$c = $ws.Range("D31")
$c.Value = "'120.70"
$c.ClearFormats()
$ws.Range("E31").Value = '  -2.69%  '
$c = $ws.Range("D32")
$c.Value = "'1.097"
$c.ClearFormats()
$ws.Range("E32").Value = '  +2.44%  '
$c = $ws.Range("D33")
$c.Value = "'5.889"
$c.ClearFormats()
$ws.Range("E33").Value = '  -2.23%  '
$c = $ws.Range("D34")
$c.Value = "'9.277"
$c.ClearFormats()
$ws.Range("E34").Value = '  -5.40%  '
$c = $ws.Range("D35")
$c.Value = "'0.08136"
$c.ClearFormats()
$ws.Range("E35").Value = '  -1.73%  '
$c = $ws.Range("D36")
$c.Value = "'1.643"
$c.ClearFormats()
$ws.Range("E36").Value = '  -16.93%  '
$c = $ws.Range("D37")
$c.Value = "'0.06207"
$c.ClearFormats()
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D38")
$c.Value = "'5.255"
$c.ClearFormats()
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D39")
$c.Value = "'0.02296"
$c.ClearFormats()
$ws.Range("E39").Value = '  -5.01%  '
$c = $ws.Range("D40")
$c.Value = "'0.2102"
$c.ClearFormats()
$ws.Range("E40").Value = '  -4.22%  '
$c = $ws.Range("D41")
$c.Value = "'1.220"
$c.ClearFormats()
$ws.Range("E41").Value = '  -4.93%  '
$c = $ws.Range("D42")
$c.Value = "'10.86"
$c.ClearFormats()
$ws.Range("E42").Value = '  -3.85%  '
$c = $ws.Range("D43")
$c.Value = "'0.9994"
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.12%  '
$c = $ws.Range("D44")
$c.Value = "'0.5931"
$c.ClearFormats()
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("E45").Value = '  -2.45%  '
$c = $ws.Range("D46")
$c.Value = "'3.721"
$c.ClearFormats()
$ws.Range("E46").Value = '  -0.85%  '
$c = $ws.Range("D47")
$c.Value = "'0.5746"
$c.ClearFormats()
$ws.Range("E47").Value = '  -3.84%  '
$c = $ws.Range("D48")
$c.Value = "'1.932"
$c.ClearFormats()
$ws.Range("E48").Value = '  -4.47%  '
$c = $ws.Range("D49")
$c.Value = "'119.42"
$c.ClearFormats()
$ws.Range("E49").Value = '  -4.19%  '
$c = $ws.Range("D50")
$c.Value = "'1.156"
$c.ClearFormats()
$ws.Range("E50").Value = '  -2.99%  '
$c = $ws.Range("D51")
$c.Value = "'0.06872"
$c.ClearFormats()
$ws.Range("E51").Value = '  -4.12%  '
